# Apply the edits described in the commit "Add files via upload / omg"
# - Update four lat/long data points on Sheet1 (rows 7 and 9)
# - Move the active cell selection from C2 to F15

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 7: long/lat
$ws.Range("B7").Value = -121.003
$ws.Range("C7").Value = 38.2965

# Row 9: long/lat
$ws.Range("B9").Value = -119.0148
$ws.Range("C9").Value = 37.5803

# Move the selection to F15 (was C2)
$ws.Range("F15").Select()
